{"js": "// Replace each unique run of text with its updated value, in document order.\n// The source document has 26 unique text runs (1 date header + 25 \"a\u00f7b=c, d\"\n// division facts), so exact, case-sensitive whole-text search/replace is\n// unambiguous and safe here.\nconst replacements = [\n  [\"2025-07-31 Thursday\", \"2025-08-01 Friday\"],\n  [\"36\u00f77=5, 1\", \"25\u00f79=2, 7\"],\n  [\"74\u00f72=37, 0\", \"39\u00f79=4, 3\"],\n  [\"19\u00f78=2, 3\", \"40\u00f79=4, 4\"],\n  [\"58\u00f76=9, 4\", \"24\u00f79=2, 6\"],\n  [\"91\u00f72=45, 1\", \"34\u00f75=6, 4\"],\n  [\"11\u00f79=1, 2\", \"51\u00f78=6, 3\"],\n  [\"29\u00f75=5, 4\", \"72\u00f76=12, 0\"],\n  [\"80\u00f74=20, 0\", \"66\u00f72=33, 0\"],\n  [\"16\u00f74=4, 0\", \"98\u00f79=10, 8\"],\n  [\"11\u00f73=3, 2\", \"95\u00f79=10, 5\"],\n  [\"86\u00f75=17, 1\", \"45\u00f77=6, 3\"],\n  [\"49\u00f74=12, 1\", \"84\u00f72=42, 0\"],\n  [\"93\u00f75=18, 3\", \"91\u00f73=30, 1\"],\n  [\"43\u00f73=14, 1\", \"85\u00f75=17, 0\"],\n  [\"47\u00f74=11, 3\", \"55\u00f79=6, 1\"],\n  [\"77\u00f76=12, 5\", \"77\u00f77=11, 0\"],\n  [\"71\u00f74=17, 3\", \"26\u00f78=3, 2\"],\n  [\"96\u00f73=32, 0\", \"88\u00f72=44, 0\"],\n  [\"92\u00f72=46, 0\", \"59\u00f77=8, 3\"],\n  [\"22\u00f79=2, 4\", \"52\u00f73=17, 1\"],\n  [\"52\u00f78=6, 4\", \"59\u00f76=9, 5\"],\n  [\"71\u00f78=8, 7\", \"27\u00f78=3, 3\"],\n  [\"87\u00f72=43, 1\", \"15\u00f76=2, 3\"],\n  [\"45\u00f75=9, 0\", \"57\u00f76=9, 3\"],\n  [\"33\u00f78=4, 1\", \"66\u00f76=11, 0\"]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  // Replace every occurrence found (there should be exactly one in this document).\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Ordered list of (old, new) text pairs -- mirrors the diff, applied as\n# whole-string, case-sensitive replacements. Every source string in this\n# document is unique, so a single Find/Replace pass per pair is unambiguous.\n$replacements = @(\n    @{old=\"2025-07-31 Thursday\"; new=\"2025-08-01 Friday\"}\n    @{old=\"36\u00f77=5, 1\"; new=\"25\u00f79=2, 7\"}\n    @{old=\"74\u00f72=37, 0\"; new=\"39\u00f79=4, 3\"}\n    @{old=\"19\u00f78=2, 3\"; new=\"40\u00f79=4, 4\"}\n    @{old=\"58\u00f76=9, 4\"; new=\"24\u00f79=2, 6\"}\n    @{old=\"91\u00f72=45, 1\"; new=\"34\u00f75=6, 4\"}\n    @{old=\"11\u00f79=1, 2\"; new=\"51\u00f78=6, 3\"}\n    @{old=\"29\u00f75=5, 4\"; new=\"72\u00f76=12, 0\"}\n    @{old=\"80\u00f74=20, 0\"; new=\"66\u00f72=33, 0\"}\n    @{old=\"16\u00f74=4, 0\"; new=\"98\u00f79=10, 8\"}\n    @{old=\"11\u00f73=3, 2\"; new=\"95\u00f79=10, 5\"}\n    @{old=\"86\u00f75=17, 1\"; new=\"45\u00f77=6, 3\"}\n    @{old=\"49\u00f74=12, 1\"; new=\"84\u00f72=42, 0\"}\n    @{old=\"93\u00f75=18, 3\"; new=\"91\u00f73=30, 1\"}\n    @{old=\"43\u00f73=14, 1\"; new=\"85\u00f75=17, 0\"}\n    @{old=\"47\u00f74=11, 3\"; new=\"55\u00f79=6, 1\"}\n    @{old=\"77\u00f76=12, 5\"; new=\"77\u00f77=11, 0\"}\n    @{old=\"71\u00f74=17, 3\"; new=\"26\u00f78=3, 2\"}\n    @{old=\"96\u00f73=32, 0\"; new=\"88\u00f72=44, 0\"}\n    @{old=\"92\u00f72=46, 0\"; new=\"59\u00f77=8, 3\"}\n    @{old=\"22\u00f79=2, 4\"; new=\"52\u00f73=17, 1\"}\n    @{old=\"52\u00f78=6, 4\"; new=\"59\u00f76=9, 5\"}\n    @{old=\"71\u00f78=8, 7\"; new=\"27\u00f78=3, 3\"}\n    @{old=\"87\u00f72=43, 1\"; new=\"15\u00f76=2, 3\"}\n    @{old=\"45\u00f75=9, 0\"; new=\"57\u00f76=9, 3\"}\n    @{old=\"33\u00f78=4, 1\"; new=\"66\u00f76=11, 0\"}\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair.old\n    $newText = $pair.new\n\n    $range = $d.Content\n    $found = $range.Find.Execute(\n        $oldText,  # FindText\n        $true,     # MatchCase\n        $false,    # MatchWholeWord\n        $false,    # MatchWildcards\n        $false,    # MatchSoundsLike\n        $false,    # MatchAllWordForms\n        $true,     # Forward\n        1,         # Wrap (wdFindContinue)\n        $false,    # Format\n        $newText,  # ReplaceWith\n        2          # Replace (wdReplaceOne)\n    )\n\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n\n"}
